# "Generate Report for Handback"
# Fills in the Latest Target File / Latest Handback File / Latest Handback
# DateTime columns for each localized-language sheet (zh-cn, de-de) now
# that the localized content has been handed back in sync with en-US, and
# flips the Status column from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$hyperlinkUnderline = 2          # xlUnderlineStyleSingle
$hyperlinkColor     = 15570276   # BGR encoding of RGB FF6495ED

function Set-HandoffHyperlink($ws, $cell, $url, $display) {
    $ws.Hyperlinks.Add($ws.Range($cell), $url, "", "", $display)
    $ws.Range($cell).Font.Underline = $hyperlinkUnderline
    $ws.Range($cell).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill Latest Target File (I), Latest Handback File (J),
#    Latest Handback DateTime (K) for both rows.
# ---------------------------------------------------------------------
Set-HandoffHyperlink $zhcn "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/addac8733cf06d7a9f8db9f69c441b1d5c5c64f4/e2e/6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md" "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md"
$zhcn.Range("J2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.27c1ee65b98a7ee1dd76788151f3087b9734550b.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-25 07:00:35"

Set-HandoffHyperlink $zhcn "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/addac8733cf06d7a9f8db9f69c441b1d5c5c64f4/e2e/e5e82743-bf71-404a-8a02-c0f5851885ee.md" "e5e82743-bf71-404a-8a02-c0f5851885ee.md"
$zhcn.Range("J3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.2fa26be12d8952b6ce93288ba391e56e600b07c7.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-25 07:00:35"

# ---------------------------------------------------------------------
# 3. de-de sheet: fill Latest Target File (I), Latest Handback File (J),
#    Latest Handback DateTime (K) for both rows.
# ---------------------------------------------------------------------
Set-HandoffHyperlink $dede "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/addac8733cf06d7a9f8db9f69c441b1d5c5c64f4/e2e/6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md" "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.md"
$dede.Range("J2").Value = "6e6e9a67-c2bc-4596-a7dc-042ca4ece98f.27c1ee65b98a7ee1dd76788151f3087b9734550b.de-de.xlf"
$dede.Range("K2").Value = "2016-08-25 07:00:43"

Set-HandoffHyperlink $dede "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/addac8733cf06d7a9f8db9f69c441b1d5c5c64f4/e2e/e5e82743-bf71-404a-8a02-c0f5851885ee.md" "e5e82743-bf71-404a-8a02-c0f5851885ee.md"
$dede.Range("J3").Value = "e5e82743-bf71-404a-8a02-c0f5851885ee.2fa26be12d8952b6ce93288ba391e56e600b07c7.de-de.xlf"
$dede.Range("K3").Value = "2016-08-25 07:00:43"

# ---------------------------------------------------------------------
# 4. Column widths widened to fit the newly-populated / longer text
#    (mirrors the auto-fit Excel performs when the report is regenerated).
#    ColumnWidth is quantized by the host in 1/6-character steps, so the
#    values below are chosen to land on (or as close as possible to) the
#    target stored widths.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.083334   # E -> stored width ~30
$overview.Columns.Item(6).ColumnWidth = 29.083334   # F -> stored width ~30

$zhcn.Columns.Item(3).ColumnWidth  = 29.083334      # C -> stored width ~30
$zhcn.Columns.Item(9).ColumnWidth  = 39.083334       # I -> stored width 40
$zhcn.Columns.Item(10).ColumnWidth = 39.083334       # J -> stored width 40

$dede.Columns.Item(3).ColumnWidth  = 29.083334      # C -> stored width ~30
$dede.Columns.Item(9).ColumnWidth  = 39.083334       # I -> stored width 40
$dede.Columns.Item(10).ColumnWidth = 39.083334       # J -> stored width 40
